# Estadisticos Segundo Parcial 23 Mayo
#
# - "Estadisticos 1P": grupo 2ARHV (row 5) gets one more aprobado
#   (a student that was reprobado is now aprobado): E 4->3, F 26->27,
#   G 86.67->90, H 7.3->7.4
# - "Estadisticos 2P": second-partial scores are filled in for every
#   grupo -- Blancos/Reprobados drop to 0, Aprobados == Totales,
#   Por_Apro == 100, and a Promedio column is populated.
# - "Estadisticos Final": recomputed the same way as 2P (all aprobados)
#   with updated Promedio values.
# - "Rescatables": the rescate list is now empty, so every student row
#   is removed, leaving only the header.

$wb = $excel.ActiveWorkbook

# --- Estadisticos 1P ---
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("E5").Value = 3
$ws1.Range("F5").Value = 27
$ws1.Range("G5").Value = 90
$ws1.Range("H5").Value = 7.4

# --- Estadisticos 2P ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
# row, Blancos, Reprobados, Aprobados, Por_Apro, Promedio
$rows2P = @(
    @(2, 0, 0, 36, 100, 6.8),
    @(3, 0, 0, 28, 100, 7.9),
    @(4, 0, 0, 23, 100, 7.3),
    @(5, 0, 0, 30, 100, 7.4),
    @(6, 0, 0, 10, 100, 10),
    @(7, 0, 0, 35, 100, 9.4),
    @(8, 0, 0, 39, 100, 8.6),
    @(9, 0, 0, 28, 100, 8.8),
    @(10, 0, 0, 25, 100, 9.6)
)
foreach ($r in $rows2P) {
    $row = $r[0]
    $ws2.Range("D$row").Value = $r[1]
    $ws2.Range("E$row").Value = $r[2]
    $ws2.Range("F$row").Value = $r[3]
    $ws2.Range("G$row").Value = $r[4]
    $ws2.Range("H$row").Value = $r[5]
}

# --- Estadisticos Final ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
# row, Reprobados, Aprobados, Por_Apro, Promedio
$rowsFinal = @(
    @(2, 0, 36, 100, 7.5),
    @(3, 0, 28, 100, 8.2),
    @(4, 0, 23, 100, 7.7),
    @(5, 0, 30, 100, 8.2),
    @(8, 0, 39, 100, 9)
)
foreach ($r in $rowsFinal) {
    $row = $r[0]
    $ws3.Range("E$row").Value = $r[1]
    $ws3.Range("F$row").Value = $r[2]
    $ws3.Range("G$row").Value = $r[3]
    $ws3.Range("H$row").Value = $r[4]
}
# rows 6, 7 and 9 already had Reprobados/Aprobados/Por_Apro at their
# final values -- only the Promedio changes
$ws3.Range("H6").Value = 9.5
$ws3.Range("H7").Value = 9.3
$ws3.Range("H9").Value = 8.7

# --- Rescatables: clear out the rescate list, keep only the header row ---
$ws4 = $wb.Worksheets.Item("Rescatables")
$ws4.Rows("2:6").Delete() | Out-Null
